$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3.0
$ws.Cells.Item(2, 7).Value2 = 57.478532
$ws.Cells.Item(2, 8).Value2 = 172.435596
$ws.Cells.Item(2, 9).Value2 = 0.2414676574042868
$ws.Cells.Item(2, 10).Value2 = 0.2414676574042868
$ws.Cells.Item(2, 11).Value2 = 3.0
$ws.Cells.Item(2, 13).Value2 = 79.465682
$ws.Cells.Item(2, 14).Value2 = 238.397046
$ws.Cells.Item(2, 15).Value2 = 0.2316778771755457
$ws.Cells.Item(2, 16).Value2 = 0.2316778771755458
$ws.Cells.Item(2, 17).Value2 = 4567.570745738824
$ws.Cells.Item(2, 18).Value2 = 41108.13671164942
$ws.Cells.Item(2, 19).Value2 = 0.05594271427397711
$ws.Cells.Item(2, 20).Value2 = 0.05594271427397712

$ws.Cells.Item(3, 5).Value2 = 3.0
$ws.Cells.Item(3, 7).Value2 = 57.478532
$ws.Cells.Item(3, 8).Value2 = 172.435596
$ws.Cells.Item(3, 9).Value2 = 0.2414676574042868
$ws.Cells.Item(3, 10).Value2 = 0.2414676574042868
$ws.Cells.Item(3, 11).Value2 = 3.0
$ws.Cells.Item(3, 13).Value2 = 21.22335033333333
$ws.Cells.Item(3, 14).Value2 = 63.670051
$ws.Cells.Item(3, 15).Value2 = 0.06187552447834749
$ws.Cells.Item(3, 16).Value2 = 0.06187552447834749
$ws.Cells.Item(3, 17).Value2 = 1219.887021281711
$ws.Cells.Item(3, 18).Value2 = 10978.9831915354
$ws.Cells.Item(3, 19).Value2 = 0.01494093794644817
$ws.Cells.Item(3, 20).Value2 = 0.01494093794644817

$ws.Cells.Item(4, 5).Value2 = 3.0
$ws.Cells.Item(4, 7).Value2 = 57.478532
$ws.Cells.Item(4, 8).Value2 = 172.435596
$ws.Cells.Item(4, 9).Value2 = 0.2414676574042868
$ws.Cells.Item(4, 10).Value2 = 0.2414676574042868
$ws.Cells.Item(4, 11).Value2 = 3.0
$ws.Cells.Item(4, 13).Value2 = 236.9265036666667
$ws.Cells.Item(4, 14).Value2 = 710.779511
$ws.Cells.Item(4, 15).Value2 = 0.6907463452728876
$ws.Cells.Item(4, 16).Value2 = 0.6907463452728876
$ws.Cells.Item(4, 17).Value2 = 13618.18762265262
$ws.Cells.Item(4, 18).Value2 = 122563.6886038736
$ws.Cells.Item(4, 19).Value2 = 0.1667929018536168
$ws.Cells.Item(4, 20).Value2 = 0.1667929018536168

$ws.Cells.Item(5, 5).Value2 = 3.0
$ws.Cells.Item(5, 7).Value2 = 57.478532
$ws.Cells.Item(5, 8).Value2 = 172.435596
$ws.Cells.Item(5, 9).Value2 = 0.2414676574042868
$ws.Cells.Item(5, 10).Value2 = 0.2414676574042868
$ws.Cells.Item(5, 11).Value2 = 3.0
$ws.Cells.Item(5, 13).Value2 = 5.385198333333332
$ws.Cells.Item(5, 14).Value2 = 16.155595
$ws.Cells.Item(5, 15).Value2 = 0.01570025307321912
$ws.Cells.Item(5, 16).Value2 = 0.01570025307321912
$ws.Cells.Item(5, 17).Value2 = 309.5332947288466
$ws.Cells.Item(5, 18).Value2 = 2785.79965255962
$ws.Cells.Item(5, 19).Value2 = 0.003791103330244675
$ws.Cells.Item(5, 20).Value2 = 0.003791103330244676

$ws.Cells.Item(6, 5).Value2 = 3.0
$ws.Cells.Item(6, 7).Value2 = 64.84043133333334
$ws.Cells.Item(6, 8).Value2 = 194.521294
$ws.Cells.Item(6, 9).Value2 = 0.272395040623924
$ws.Cells.Item(6, 10).Value2 = 0.2723950406239241
$ws.Cells.Item(6, 11).Value2 = 3.0
$ws.Cells.Item(6, 13).Value2 = 79.465682
$ws.Cells.Item(6, 14).Value2 = 238.397046
$ws.Cells.Item(6, 15).Value2 = 0.2316778771755457
$ws.Cells.Item(6, 16).Value2 = 0.2316778771755458
$ws.Cells.Item(6, 17).Value2 = 5152.589097077504
$ws.Cells.Item(6, 18).Value2 = 46373.30187369753
$ws.Cells.Item(6, 19).Value2 = 0.06310790476489726
$ws.Cells.Item(6, 20).Value2 = 0.06310790476489728

$ws.Cells.Item(7, 5).Value2 = 3.0
$ws.Cells.Item(7, 7).Value2 = 64.84043133333334
$ws.Cells.Item(7, 8).Value2 = 194.521294
$ws.Cells.Item(7, 9).Value2 = 0.272395040623924
$ws.Cells.Item(7, 10).Value2 = 0.2723950406239241
$ws.Cells.Item(7, 11).Value2 = 3.0
$ws.Cells.Item(7, 13).Value2 = 21.22335033333333
$ws.Cells.Item(7, 14).Value2 = 63.670051
$ws.Cells.Item(7, 15).Value2 = 0.06187552447834749
$ws.Cells.Item(7, 16).Value2 = 0.06187552447834749
$ws.Cells.Item(7, 17).Value2 = 1376.131189951777
$ws.Cells.Item(7, 18).Value2 = 12385.180709566
$ws.Cells.Item(7, 19).Value2 = 0.01685458600390607
$ws.Cells.Item(7, 20).Value2 = 0.01685458600390607

$ws.Cells.Item(8, 5).Value2 = 3.0
$ws.Cells.Item(8, 7).Value2 = 64.84043133333334
$ws.Cells.Item(8, 8).Value2 = 194.521294
$ws.Cells.Item(8, 9).Value2 = 0.272395040623924
$ws.Cells.Item(8, 10).Value2 = 0.2723950406239241
$ws.Cells.Item(8, 11).Value2 = 3.0
$ws.Cells.Item(8, 13).Value2 = 236.9265036666667
$ws.Cells.Item(8, 14).Value2 = 710.779511
$ws.Cells.Item(8, 15).Value2 = 0.6907463452728876
$ws.Cells.Item(8, 16).Value2 = 0.6907463452728876
$ws.Cells.Item(8, 17).Value2 = 15362.41669204525
$ws.Cells.Item(8, 18).Value2 = 138261.7502284072
$ws.Cells.Item(8, 19).Value2 = 0.1881558787814352
$ws.Cells.Item(8, 20).Value2 = 0.1881558787814353

$ws.Cells.Item(9, 5).Value2 = 3.0
$ws.Cells.Item(9, 7).Value2 = 64.84043133333334
$ws.Cells.Item(9, 8).Value2 = 194.521294
$ws.Cells.Item(9, 9).Value2 = 0.272395040623924
$ws.Cells.Item(9, 10).Value2 = 0.2723950406239241
$ws.Cells.Item(9, 11).Value2 = 3.0
$ws.Cells.Item(9, 13).Value2 = 5.385198333333332
$ws.Cells.Item(9, 14).Value2 = 16.155595
$ws.Cells.Item(9, 15).Value2 = 0.01570025307321912
$ws.Cells.Item(9, 16).Value2 = 0.01570025307321912
$ws.Cells.Item(9, 17).Value2 = 349.1785827488811
$ws.Cells.Item(9, 18).Value2 = 3142.60724473993
$ws.Cells.Item(9, 19).Value2 = 0.00427667107368541
$ws.Cells.Item(9, 20).Value2 = 0.004276671073685411

$ws.Cells.Item(10, 5).Value2 = 3.0
$ws.Cells.Item(10, 7).Value2 = 85.31555666666667
$ws.Cells.Item(10, 8).Value2 = 255.94667
$ws.Cells.Item(10, 9).Value2 = 0.3584111648579104
$ws.Cells.Item(10, 10).Value2 = 0.3584111648579105
$ws.Cells.Item(10, 11).Value2 = 3.0
$ws.Cells.Item(10, 13).Value2 = 79.465682
$ws.Cells.Item(10, 14).Value2 = 238.397046
$ws.Cells.Item(10, 15).Value2 = 0.2316778771755457
$ws.Cells.Item(10, 16).Value2 = 0.2316778771755458
$ws.Cells.Item(10, 17).Value2 = 6779.658895726313
$ws.Cells.Item(10, 18).Value2 = 61016.93006153682
$ws.Cells.Item(10, 19).Value2 = 0.08303593783029525
$ws.Cells.Item(10, 20).Value2 = 0.08303593783029527

$ws.Cells.Item(11, 5).Value2 = 3.0
$ws.Cells.Item(11, 7).Value2 = 85.31555666666667
$ws.Cells.Item(11, 8).Value2 = 255.94667
$ws.Cells.Item(11, 9).Value2 = 0.3584111648579104
$ws.Cells.Item(11, 10).Value2 = 0.3584111648579105
$ws.Cells.Item(11, 11).Value2 = 3.0
$ws.Cells.Item(11, 13).Value2 = 21.22335033333333
$ws.Cells.Item(11, 14).Value2 = 63.670051
$ws.Cells.Item(11, 15).Value2 = 0.06187552447834749
$ws.Cells.Item(11, 16).Value2 = 0.06187552447834749
$ws.Cells.Item(11, 17).Value2 = 1810.681948020019
$ws.Cells.Item(11, 18).Value2 = 16296.13753218017
$ws.Cells.Item(11, 19).Value2 = 0.02217687880447867
$ws.Cells.Item(11, 20).Value2 = 0.02217687880447868

$ws.Cells.Item(12, 5).Value2 = 3.0
$ws.Cells.Item(12, 7).Value2 = 85.31555666666667
$ws.Cells.Item(12, 8).Value2 = 255.94667
$ws.Cells.Item(12, 9).Value2 = 0.3584111648579104
$ws.Cells.Item(12, 10).Value2 = 0.3584111648579105
$ws.Cells.Item(12, 11).Value2 = 3.0
$ws.Cells.Item(12, 13).Value2 = 236.9265036666667
$ws.Cells.Item(12, 14).Value2 = 710.779511
$ws.Cells.Item(12, 15).Value2 = 0.6907463452728876
$ws.Cells.Item(12, 16).Value2 = 0.6907463452728876
$ws.Cells.Item(12, 17).Value2 = 20213.51654940871
$ws.Cells.Item(12, 18).Value2 = 181921.6489446783
$ws.Cells.Item(12, 19).Value2 = 0.2475712022306
$ws.Cells.Item(12, 20).Value2 = 0.2475712022306001

$ws.Cells.Item(13, 5).Value2 = 3.0
$ws.Cells.Item(13, 7).Value2 = 85.31555666666667
$ws.Cells.Item(13, 8).Value2 = 255.94667
$ws.Cells.Item(13, 9).Value2 = 0.3584111648579104
$ws.Cells.Item(13, 10).Value2 = 0.3584111648579105
$ws.Cells.Item(13, 11).Value2 = 3.0
$ws.Cells.Item(13, 13).Value2 = 5.385198333333332
$ws.Cells.Item(13, 14).Value2 = 16.155595
$ws.Cells.Item(13, 15).Value2 = 0.01570025307321912
$ws.Cells.Item(13, 16).Value2 = 0.01570025307321912
$ws.Cells.Item(13, 17).Value2 = 459.4411935687388
$ws.Cells.Item(13, 18).Value2 = 4134.970742118649
$ws.Cells.Item(13, 19).Value2 = 0.005627145992536453
$ws.Cells.Item(13, 20).Value2 = 0.005627145992536454

$ws.Cells.Item(14, 5).Value2 = 3.0
$ws.Cells.Item(14, 7).Value2 = 30.40370266666666
$ws.Cells.Item(14, 8).Value2 = 91.211108
$ws.Cells.Item(14, 9).Value2 = 0.1277261371138787
$ws.Cells.Item(14, 10).Value2 = 0.1277261371138788
$ws.Cells.Item(14, 11).Value2 = 3.0
$ws.Cells.Item(14, 13).Value2 = 79.465682
$ws.Cells.Item(14, 14).Value2 = 238.397046
$ws.Cells.Item(14, 15).Value2 = 0.2316778771755457
$ws.Cells.Item(14, 16).Value2 = 0.2316778771755458
$ws.Cells.Item(14, 17).Value2 = 2416.050967731885
$ws.Cells.Item(14, 18).Value2 = 21744.45870958697
$ws.Cells.Item(14, 19).Value2 = 0.02959132030637611
$ws.Cells.Item(14, 20).Value2 = 0.02959132030637612

$ws.Cells.Item(15, 5).Value2 = 3.0
$ws.Cells.Item(15, 7).Value2 = 30.40370266666666
$ws.Cells.Item(15, 8).Value2 = 91.211108
$ws.Cells.Item(15, 9).Value2 = 0.1277261371138787
$ws.Cells.Item(15, 10).Value2 = 0.1277261371138788
$ws.Cells.Item(15, 11).Value2 = 3.0
$ws.Cells.Item(15, 13).Value2 = 21.22335033333333
$ws.Cells.Item(15, 14).Value2 = 63.670051
$ws.Cells.Item(15, 15).Value2 = 0.06187552447834749
$ws.Cells.Item(15, 16).Value2 = 0.06187552447834749
$ws.Cells.Item(15, 17).Value2 = 645.2684331251675
$ws.Cells.Item(15, 18).Value2 = 5807.415898126508
$ws.Cells.Item(15, 19).Value2 = 0.00790312172351457
$ws.Cells.Item(15, 20).Value2 = 0.007903121723514574

$ws.Cells.Item(16, 5).Value2 = 3.0
$ws.Cells.Item(16, 7).Value2 = 30.40370266666666
$ws.Cells.Item(16, 8).Value2 = 91.211108
$ws.Cells.Item(16, 9).Value2 = 0.1277261371138787
$ws.Cells.Item(16, 10).Value2 = 0.1277261371138788
$ws.Cells.Item(16, 11).Value2 = 3.0
$ws.Cells.Item(16, 13).Value2 = 236.9265036666667
$ws.Cells.Item(16, 14).Value2 = 710.779511
$ws.Cells.Item(16, 15).Value2 = 0.6907463452728876
$ws.Cells.Item(16, 16).Value2 = 0.6907463452728876
$ws.Cells.Item(16, 17).Value2 = 7203.442971334242
$ws.Cells.Item(16, 18).Value2 = 64830.98674200818
$ws.Cells.Item(16, 19).Value2 = 0.08822636240723546
$ws.Cells.Item(16, 20).Value2 = 0.08822636240723548

$ws.Cells.Item(17, 5).Value2 = 3.0
$ws.Cells.Item(17, 7).Value2 = 30.40370266666666
$ws.Cells.Item(17, 8).Value2 = 91.211108
$ws.Cells.Item(17, 9).Value2 = 0.1277261371138787
$ws.Cells.Item(17, 10).Value2 = 0.1277261371138788
$ws.Cells.Item(17, 11).Value2 = 3.0
$ws.Cells.Item(17, 13).Value2 = 5.385198333333332
$ws.Cells.Item(17, 14).Value2 = 16.155595
$ws.Cells.Item(17, 15).Value2 = 0.01570025307321912
$ws.Cells.Item(17, 16).Value2 = 0.01570025307321912
$ws.Cells.Item(17, 17).Value2 = 163.7299689276945
$ws.Cells.Item(17, 18).Value2 = 1473.56972034926
$ws.Cells.Item(17, 19).Value2 = 0.002005332676752581
$ws.Cells.Item(17, 20).Value2 = 0.002005332676752582
